$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$hf = $m.HeadersFooters
$dt = $hf.DateAndTime
Write-Output "Initial Text: $($dt.Text)"
Write-Output "Initial UseFormat: $($dt.UseFormat)"
$dt.UseFormat = 0
$dt.Text = "2/17/2018"
Write-Output "After set (UseFormat=0): $($dt.Text)"
$dt.UseFormat = -1
Write-Output "After UseFormat=-1: $($dt.Text)"
$sh = $m.Shapes.Item(3)
Write-Output "Shape Text: [$($sh.TextFrame.TextRange.Text)]"
